$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.606.89'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.27'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6982'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07711'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3061'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.66'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07750'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.164'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('D13').Value = '1.855.15'
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.34'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.37%  '
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.569'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '29.595.75'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '2.106.21'
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.0000'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.615'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.98%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1504'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.915'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.40'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.535'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.179'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.193'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05090'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7763'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.895'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.153'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.686'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('D38').Value = '1.324.33'
$ws.Range('E38').Value = '  +10.20%  '
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.733'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9595'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.44'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.827'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +11.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.781'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000125'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').Value = '2.005.10'
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5218'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.57'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.92%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.780'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.37%  '
$ws.Range('E51').Value = '  +1.55%  '
